$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 360, shifting existing rows 360-479 down to 361-480
$ws.Rows.Item(360).Insert()

# Populate the newly inserted row 360 with the new weekly price-report entry
$ws.Range("A360").Value = 7
$ws.Range("B360").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C360").Value = "Ñuble"
$ws.Range("D360").Value = 44524
$ws.Range("E360").Value = 16
$ws.Range("F360").Value = "Fruta"
$ws.Range("G360").Value = 100102
$ws.Range("H360").Value = "Cítricos"
$ws.Range("I360").Value = 100102003
$ws.Range("J360").Value = "Limón"
$ws.Range("K360").Value = "Sin especificar"
$ws.Range("L360").Value = "1a amarillo"
$ws.Range("M360").Value = 160
$ws.Range("N360").Value = 7000
$ws.Range("O360").Value = 7500
$ws.Range("P360").Value = 7250
$ws.Range("Q360").Value = "$/malla 16 kilos"
$ws.Range("R360").Value = "Región de O'Higgins"
$ws.Range("S360").Value = 453
$ws.Range("T360").Value = 16
